$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.539.96"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "3.593.15"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'607.54"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'149.10"
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("D7").Value = "3.591.85"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'0.136"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "'7.97"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").Value = "'0.414"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "4.198.01"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "'0.0000206"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "'29.69"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "3.587.04"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").Value = "66.545.74"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "'11.07"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "'6.34"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").Value = "'14.88"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "'424.82"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "'0.614"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'78.19"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("E28").Value = "  +3.32%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D31").Value = "3.587.29"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  +4.51%  "
$ws.Range("D33").Value = "'25.05"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D36").Value = "'7.75"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").Value = "'5.55"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("D39").Value = "'175.01"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'0.0854"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'0.882"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +4.11%  "
$ws.Range("D47").Value = "'23.76"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").Value = "'24.39"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("E49").Value = "  -4.13%  "
$ws.Range("D50").Value = "'7.14"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'0.942"
$ws.Range("E51").Value = "  +0.91%  "
